$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.147.89'
$ws.Range('E2').Value = '  +1.38%  '
$ws.Range('D3').Value = '3.907.28'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '464.89'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +8.30%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.80'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.60%  '
$ws.Range('E7').Value = '  +2.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('E9').Value = '  -1.48%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.164'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +5.26%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0000339'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +4.70%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '42.98'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.78%  '
$ws.Range('B13').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C13').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D13').Value = '4.528.50'
$ws.Range('E13').Value = '  +3.16%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '10.37'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.03%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.28'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.77%  '
$ws.Range('D16').Value = '3.889.10'
$ws.Range('E16').Value = '  +2.86%  '
$ws.Range('E17').Value = '  -0.50%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '19.97'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.89%  '
$ws.Range('E19').Value = '  +1.74%  '
$ws.Range('D20').Value = '67.333.15'
$ws.Range('E20').Value = '  +1.43%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '431.55'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +5.17%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '14.67'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -4.23%  '
$ws.Range('E23').Value = '  +1.00%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '89.05'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +4.10%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '38.66'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +4.25%  '
$ws.Range('E26').Value = '  +6.44%  '
$ws.Range('E27').Value = '  +5.18%  '
$ws.Range('E28').Value = '  +2.51%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.60'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.58%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '739.78'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +5.72%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '13.61'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.94%  '
$ws.Range('E32').Value = '  +1.81%  '
$ws.Range('E33').Value = '  -0.12%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '43.07'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +5.60%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.157'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +3.95%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '58.16'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.65%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.00'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.09%  '
$ws.Range('D38').Value = '0.0₃0783'
$ws.Range('E38').Value = '  +15.06%  '
$ws.Range('E39').Value = '  -6.27%  '
$ws.Range('E40').Value = '  +11.49%  '
$ws.Range('E41').Value = '  -0.43%  '
$ws.Range('E42').Value = '  -1.04%  '
$ws.Range('E43').Value = '  +0.08%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.334'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +4.00%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.78'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +4.62%  '
$ws.Range('E46').Value = '  +4.49%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.41'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.41%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.49'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.55%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.14'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.82%  '
$ws.Range('E50').Value = '  +2.19%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '143.46'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.44%  '
